# Support LINK price feed from USD to wei
# OKLGAtomicSwapInstance (row 9) base USD price changes from 0 to 2.
# Excel automatically recalculates the dependent per-chain price formulas
# in D9:K9 (each is $C9/<chain base price in row 2>).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("C9").Value = 2
